$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.138.57'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.674.66'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '214.10'
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '22.83'
$ws.Range('E8').Value = '  +6.69%  '
$ws.Range('E9').Value = '  +2.97%  '
$ws.Range('D10').Value = '0.0620'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '1.912.15'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.682.97'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('D14').Value = '4.20'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '0.557'
$ws.Range('E15').Value = '  +4.40%  '
$ws.Range('D16').Value = '66.47'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = '27.092.83'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '234.72'
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0741'
$ws.Range('E19').Value = '  +0.42%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '7.82'
$ws.Range('E20').Value = '  -4.13%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '4.53'
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('E23').Value = '  +3.01%  '
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('D25').Value = '147.86'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').Value = '7.45'
$ws.Range('E26').Value = '  +2.33%  '
$ws.Range('D27').Value = '16.37'
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = '0.0499'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = '1.540.62'
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').Value = '3.22'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').Value = '1.65'
$ws.Range('E35').Value = '  -4.07%  '
$ws.Range('E36').Value = '  +3.15%  '
$ws.Range('D37').Value = '0.945'
$ws.Range('E37').Value = '  +3.17%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('E40').Value = '  +2.30%  '
$ws.Range('D41').Value = '69.91'
$ws.Range('E41').Value = '  +3.02%  '
$ws.Range('D42').Value = '5.77'
$ws.Range('E42').Value = '  +4.35%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = '1.820.79'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('D47').Value = '89.49'
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('E48').Value = '  +6.25%  '
$ws.Range('D49').Value = '0.0₆0110'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('D50').Value = '8.18'
$ws.Range('E50').Value = '  +1.53%  '
$ws.Range('E51').Value = '  -0.10%  '
